# "added org update and delete test case"
# The test data's second row held an org name placeholder ("AUTO_ORG_l").
# Replace it with a new distinct org name placeholder ("AUTO_ORG_ERZYN")
# so the new update/delete test case has its own org to exercise.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "AUTO_ORG_ERZYN"
